$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (written first so these become shared strings 0..3) ---
$ws.Range("A1").Value = "Reference No"
$ws.Range("B1").Value = "Customer"
$ws.Range("C1").Value = "Created At"
$ws.Range("D1").Value = "Service Cost"

# --- Column A: reference numbers (shared strings 4..6) ---
$ws.Range("A2").Value = "REF-20250731-001"
$ws.Range("A3").Value = "REF-20250731-002"
$ws.Range("A4").Value = "REF-20250731-003"

# --- Column B: customer names (shared strings 7..9) ---
$ws.Range("B2").Value = "CRESTANKS LIMITED"
$ws.Range("B3").Value = "KYAGALANYI COFFEE LIMITED"
$ws.Range("B4").Value = "UGANDA REVENUE AUTHORITY"

# --- Column C: date-looking text that must stay literal text, not get
# auto-converted to a date serial number, and must not leave behind any
# unused style/numFmt entries. Writing it as a formula that yields the
# string, then collapsing it to a static value via Copy/PasteSpecial
# (values only) avoids Excel's text-to-date autodetection entirely, and
# does not touch the style tables. (shared string 10, first use in C2) ---
foreach ($row in 2..4) {
    $cell = $ws.Range("C$row")
    $cell.Formula = '="2025-07-30"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

# --- Column D: plain numbers ---
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 10

# --- Header formatting: bold, thin box border, centered horizontal,
# top vertical. Build the whole style on a single cell (A1) so all the
# property writes coalesce into exactly one new cellXf, then propagate
# that exact style to B1:D1 via copy/paste-format so no extra cellXfs
# get allocated. ---
$hdr = $ws.Range("A1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
